$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before the current row 31, shifting existing rows 31-150 down to 32-151.
$ws.Rows(31).Insert()

# Populate the newly inserted row 31 with the new weekly record.
$ws.Range("A31").Value = 11
$ws.Range("B31").Value = "Vega Monumental Concepción"
$ws.Range("C31").Value = "Bíobío"
$ws.Range("D31").Value = 45037
$ws.Range("E31").Value = 8
$ws.Range("F31").Value = 100112001
$ws.Range("G31").Value = "Berenjena"
$ws.Range("H31").Value = "Sin especificar"
$ws.Range("I31").Value = "Primera"
$ws.Range("J31").Value = 100
$ws.Range("K31").Value = 7500
$ws.Range("L31").Value = 8000
$ws.Range("M31").Value = 7750
$ws.Range("N31").Value = '$/caja 60 unidades'
$ws.Range("O31").Value = "Región de Arica y Parinacota"
$ws.Range("P31").Value = 129
$ws.Range("Q31").Value = 60
$ws.Range("R31").Value = "Hortaliza"
